# Update Ncam1-Fgfr1 LR-pair sheet with recomputed TPM values and the new
# "Resolving-Mac" cluster (adds 6 rows: Resolving-Mac as sending cluster,
# and a Resolving-Mac column for every existing sending cluster).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = New-Object 'object[,]' 36,20
# row 2: ECs -> ECs
$data[0,0] = "ECs"
$data[0,1] = "Ncam1"
$data[0,2] = "Fgfr1"
$data[0,3] = "ECs"
$data[0,4] = 3
$data[0,5] = 1
$data[0,6] = 0.9949870000000001
$data[0,7] = 2.984961
$data[0,8] = 0.03855738270564991
$data[0,9] = 0.03855738270564991
$data[0,10] = 3
$data[0,11] = 1
$data[0,12] = 10.48767733333333
$data[0,13] = 31.463032
$data[0,14] = 0.1222087640673552
$data[0,15] = 0.1222087640673552
$data[0,16] = 10.43510260686133
$data[0,17] = 93.915923461752
$data[0,18] = 0.004712050086129493
$data[0,19] = 0.004712050086129493
# row 3: ECs -> FAPs
$data[1,0] = "ECs"
$data[1,1] = "Ncam1"
$data[1,2] = "Fgfr1"
$data[1,3] = "FAPs"
$data[1,4] = 3
$data[1,5] = 1
$data[1,6] = 0.9949870000000001
$data[1,7] = 2.984961
$data[1,8] = 0.03855738270564991
$data[1,9] = 0.03855738270564991
$data[1,10] = 3
$data[1,11] = 1
$data[1,12] = 62.99699166666667
$data[1,13] = 188.990975
$data[1,14] = 0.7340790765058636
$data[1,15] = 0.7340790765058635
$data[1,16] = 62.68118774744167
$data[1,17] = 564.130689726975
$data[1,18] = 0.02830416788904664
$data[1,19] = 0.02830416788904664
# row 4: ECs -> Inflammatory-Mac
$data[2,0] = "ECs"
$data[2,1] = "Ncam1"
$data[2,2] = "Fgfr1"
$data[2,3] = "Inflammatory-Mac"
$data[2,4] = 3
$data[2,5] = 1
$data[2,6] = 0.9949870000000001
$data[2,7] = 2.984961
$data[2,8] = 0.03855738270564991
$data[2,9] = 0.03855738270564991
$data[2,10] = 3
$data[2,11] = 1
$data[2,12] = 0.3322793333333333
$data[2,13] = 0.996838
$data[2,14] = 0.003871919907635547
$data[2,15] = 0.003871919907635547
$data[2,16] = 0.3306136170353333
$data[2,17] = 2.975522553318
$data[2,18] = 0.0001492910976843284
$data[2,19] = 0.0001492910976843285
# row 5: ECs -> MuSCs
$data[3,0] = "ECs"
$data[3,1] = "Ncam1"
$data[3,2] = "Fgfr1"
$data[3,3] = "MuSCs"
$data[3,4] = 3
$data[3,5] = 1
$data[3,6] = 0.9949870000000001
$data[3,7] = 2.984961
$data[3,8] = 0.03855738270564991
$data[3,9] = 0.03855738270564991
$data[3,10] = 3
$data[3,11] = 1
$data[3,12] = 10.25458433333333
$data[3,13] = 30.763753
$data[3,14] = 0.1194926233493133
$data[3,15] = 0.1194926233493133
$data[3,16] = 10.20317810207033
$data[3,17] = 91.82860291863301
$data[3,18] = 0.004607322808981552
$data[3,19] = 0.004607322808981552
# row 6: ECs -> Neutrophils
$data[4,0] = "ECs"
$data[4,1] = "Ncam1"
$data[4,2] = "Fgfr1"
$data[4,3] = "Neutrophils"
$data[4,4] = 3
$data[4,5] = 1
$data[4,6] = 0.9949870000000001
$data[4,7] = 2.984961
$data[4,8] = 0.03855738270564991
$data[4,9] = 0.03855738270564991
$data[4,10] = 3
$data[4,11] = 1
$data[4,12] = 0.7572163333333334
$data[4,13] = 2.271649
$data[4,14] = 0.008823543029319092
$data[4,15] = 0.00882354302931909
$data[4,16] = 0.7534204078543334
$data[4,17] = 6.780783670689001
$data[4,18] = 0.0003402127254012257
$data[4,19] = 0.0003402127254012257
# row 7: ECs -> Resolving-Mac
$data[5,0] = "ECs"
$data[5,1] = "Ncam1"
$data[5,2] = "Fgfr1"
$data[5,3] = "Resolving-Mac"
$data[5,4] = 3
$data[5,5] = 1
$data[5,6] = 0.9949870000000001
$data[5,7] = 2.984961
$data[5,8] = 0.03855738270564991
$data[5,9] = 0.03855738270564991
$data[5,10] = 3
$data[5,11] = 1
$data[5,12] = 0.9889696666666666
$data[5,13] = 2.966909
$data[5,14] = 0.01152407314051338
$data[5,15] = 0.01152407314051338
$data[5,16] = 0.9840119617276667
$data[5,17] = 8.856107655549
$data[5,18] = 0.0004443380984066751
$data[5,19] = 0.0004443380984066751
# row 8: FAPs -> ECs
$data[6,0] = "FAPs"
$data[6,1] = "Ncam1"
$data[6,2] = "Fgfr1"
$data[6,3] = "ECs"
$data[6,4] = 3
$data[6,5] = 1
$data[6,6] = 8.469728666666667
$data[6,7] = 25.409186
$data[6,8] = 0.3282159160005915
$data[6,9] = 0.3282159160005916
$data[6,10] = 3
$data[6,11] = 1
$data[6,12] = 10.48767733333333
$data[6,13] = 31.463032
$data[6,14] = 0.1222087640673552
$data[6,15] = 0.1222087640673552
$data[6,16] = 88.82778135688355
$data[6,17] = 799.450032211952
$data[6,18] = 0.04011086144166718
$data[6,19] = 0.04011086144166719
# row 9: FAPs -> FAPs
$data[7,0] = "FAPs"
$data[7,1] = "Ncam1"
$data[7,2] = "Fgfr1"
$data[7,3] = "FAPs"
$data[7,4] = 3
$data[7,5] = 1
$data[7,6] = 8.469728666666667
$data[7,7] = 25.409186
$data[7,8] = 0.3282159160005915
$data[7,9] = 0.3282159160005916
$data[7,10] = 3
$data[7,11] = 1
$data[7,12] = 62.99699166666667
$data[7,13] = 188.990975
$data[7,14] = 0.7340790765058636
$data[7,15] = 0.7340790765058635
$data[7,16] = 533.5674262329278
$data[7,17] = 4802.10683609635
$data[7,18] = 0.2409364365122403
$data[7,19] = 0.2409364365122404
# row 10: FAPs -> Inflammatory-Mac
$data[8,0] = "FAPs"
$data[8,1] = "Ncam1"
$data[8,2] = "Fgfr1"
$data[8,3] = "Inflammatory-Mac"
$data[8,4] = 3
$data[8,5] = 1
$data[8,6] = 8.469728666666667
$data[8,7] = 25.409186
$data[8,8] = 0.3282159160005915
$data[8,9] = 0.3282159160005916
$data[8,10] = 3
$data[8,11] = 1
$data[8,12] = 0.3322793333333333
$data[8,13] = 0.996838
$data[8,14] = 0.003871919907635547
$data[8,15] = 0.003871919907635547
$data[8,16] = 2.814315794874222
$data[8,17] = 25.328842153868
$data[8,18] = 0.001270825739165527
$data[8,19] = 0.001270825739165527
# row 11: FAPs -> MuSCs
$data[9,0] = "FAPs"
$data[9,1] = "Ncam1"
$data[9,2] = "Fgfr1"
$data[9,3] = "MuSCs"
$data[9,4] = 3
$data[9,5] = 1
$data[9,6] = 8.469728666666667
$data[9,7] = 25.409186
$data[9,8] = 0.3282159160005915
$data[9,9] = 0.3282159160005916
$data[9,10] = 3
$data[9,11] = 1
$data[9,12] = 10.25458433333333
$data[9,13] = 30.763753
$data[9,14] = 0.1194926233493133
$data[9,15] = 0.1194926233493133
$data[9,16] = 86.85354689278422
$data[9,17] = 781.681922035058
$data[9,18] = 0.03921938082790855
$data[9,19] = 0.03921938082790855
# row 12: FAPs -> Neutrophils
$data[10,0] = "FAPs"
$data[10,1] = "Ncam1"
$data[10,2] = "Fgfr1"
$data[10,3] = "Neutrophils"
$data[10,4] = 3
$data[10,5] = 1
$data[10,6] = 8.469728666666667
$data[10,7] = 25.409186
$data[10,8] = 0.3282159160005915
$data[10,9] = 0.3282159160005916
$data[10,10] = 3
$data[10,11] = 1
$data[10,12] = 0.7572163333333334
$data[10,13] = 2.271649
$data[10,14] = 0.008823543029319092
$data[10,15] = 0.00882354302931909
$data[10,16] = 6.413416885301556
$data[10,17] = 57.72075196771401
$data[10,18] = 0.0028960272577386
$data[10,19] = 0.0028960272577386
# row 13: FAPs -> Resolving-Mac
$data[11,0] = "FAPs"
$data[11,1] = "Ncam1"
$data[11,2] = "Fgfr1"
$data[11,3] = "Resolving-Mac"
$data[11,4] = 3
$data[11,5] = 1
$data[11,6] = 8.469728666666667
$data[11,7] = 25.409186
$data[11,8] = 0.3282159160005915
$data[11,9] = 0.3282159160005916
$data[11,10] = 3
$data[11,11] = 1
$data[11,12] = 0.9889696666666666
$data[11,13] = 2.966909
$data[11,14] = 0.01152407314051338
$data[11,15] = 0.01152407314051338
$data[11,16] = 8.376304736230445
$data[11,17] = 75.386742626074
$data[11,18] = 0.003782384221871412
$data[11,19] = 0.003782384221871412
# row 14: Inflammatory-Mac -> ECs
$data[12,0] = "Inflammatory-Mac"
$data[12,1] = "Ncam1"
$data[12,2] = "Fgfr1"
$data[12,3] = "ECs"
$data[12,4] = 2
$data[12,5] = 0.6666666666666666
$data[12,6] = 0.04495399999999999
$data[12,7] = 0.134862
$data[12,8] = 0.001742041435867791
$data[12,9] = 0.001742041435867791
$data[12,10] = 3
$data[12,11] = 1
$data[12,12] = 10.48767733333333
$data[12,13] = 31.463032
$data[12,14] = 0.1222087640673552
$data[12,15] = 0.1222087640673552
$data[12,16] = 0.4714630468426665
$data[12,17] = 4.243167421583999
$data[12,18] = 0.0002128927308315236
$data[12,19] = 0.0002128927308315236
# row 15: Inflammatory-Mac -> FAPs
$data[13,0] = "Inflammatory-Mac"
$data[13,1] = "Ncam1"
$data[13,2] = "Fgfr1"
$data[13,3] = "FAPs"
$data[13,4] = 2
$data[13,5] = 0.6666666666666666
$data[13,6] = 0.04495399999999999
$data[13,7] = 0.134862
$data[13,8] = 0.001742041435867791
$data[13,9] = 0.001742041435867791
$data[13,10] = 3
$data[13,11] = 1
$data[13,12] = 62.99699166666667
$data[13,13] = 188.990975
$data[13,14] = 0.7340790765058636
$data[13,15] = 0.7340790765058635
$data[13,16] = 2.831966763383333
$data[13,17] = 25.48770087044999
$data[13,18] = 0.001278796168476776
$data[13,19] = 0.001278796168476777
# row 16: Inflammatory-Mac -> Inflammatory-Mac
$data[14,0] = "Inflammatory-Mac"
$data[14,1] = "Ncam1"
$data[14,2] = "Fgfr1"
$data[14,3] = "Inflammatory-Mac"
$data[14,4] = 2
$data[14,5] = 0.6666666666666666
$data[14,6] = 0.04495399999999999
$data[14,7] = 0.134862
$data[14,8] = 0.001742041435867791
$data[14,9] = 0.001742041435867791
$data[14,10] = 3
$data[14,11] = 1
$data[14,12] = 0.3322793333333333
$data[14,13] = 0.996838
$data[14,14] = 0.003871919907635547
$data[14,15] = 0.003871919907635547
$data[14,16] = 0.01493728515066666
$data[14,17] = 0.134435566356
$data[14,18] = 0.000006745044915462512
$data[14,19] = 0.000006745044915462514
# row 17: Inflammatory-Mac -> MuSCs
$data[15,0] = "Inflammatory-Mac"
$data[15,1] = "Ncam1"
$data[15,2] = "Fgfr1"
$data[15,3] = "MuSCs"
$data[15,4] = 2
$data[15,5] = 0.6666666666666666
$data[15,6] = 0.04495399999999999
$data[15,7] = 0.134862
$data[15,8] = 0.001742041435867791
$data[15,9] = 0.001742041435867791
$data[15,10] = 3
$data[15,11] = 1
$data[15,12] = 10.25458433333333
$data[15,13] = 30.763753
$data[15,14] = 0.1194926233493133
$data[15,15] = 0.1194926233493133
$data[15,16] = 0.4609845841206666
$data[15,17] = 4.148861257086
$data[15,18] = 0.0002081611011550469
$data[15,19] = 0.0002081611011550469
# row 18: Inflammatory-Mac -> Neutrophils
$data[16,0] = "Inflammatory-Mac"
$data[16,1] = "Ncam1"
$data[16,2] = "Fgfr1"
$data[16,3] = "Neutrophils"
$data[16,4] = 2
$data[16,5] = 0.6666666666666666
$data[16,6] = 0.04495399999999999
$data[16,7] = 0.134862
$data[16,8] = 0.001742041435867791
$data[16,9] = 0.001742041435867791
$data[16,10] = 3
$data[16,11] = 1
$data[16,12] = 0.7572163333333334
$data[16,13] = 2.271649
$data[16,14] = 0.008823543029319092
$data[16,15] = 0.00882354302931909
$data[16,16] = 0.03403990304866666
$data[16,17] = 0.306359127438
$data[16,18] = 0.00001537097756823627
$data[16,19] = 0.00001537097756823627
# row 19: Inflammatory-Mac -> Resolving-Mac
$data[17,0] = "Inflammatory-Mac"
$data[17,1] = "Ncam1"
$data[17,2] = "Fgfr1"
$data[17,3] = "Resolving-Mac"
$data[17,4] = 2
$data[17,5] = 0.6666666666666666
$data[17,6] = 0.04495399999999999
$data[17,7] = 0.134862
$data[17,8] = 0.001742041435867791
$data[17,9] = 0.001742041435867791
$data[17,10] = 3
$data[17,11] = 1
$data[17,12] = 0.9889696666666666
$data[17,13] = 2.966909
$data[17,14] = 0.01152407314051338
$data[17,15] = 0.01152407314051338
$data[17,16] = 0.04445814239533333
$data[17,17] = 0.4001232815579999
$data[17,18] = 0.00002007541292074537
$data[17,19] = 0.00002007541292074537
# row 20: MuSCs -> ECs
$data[18,0] = "MuSCs"
$data[18,1] = "Ncam1"
$data[18,2] = "Fgfr1"
$data[18,3] = "ECs"
$data[18,4] = 3
$data[18,5] = 1
$data[18,6] = 15.96019966666667
$data[18,7] = 47.880599
$data[18,8] = 0.6184839868322428
$data[18,9] = 0.6184839868322429
$data[18,10] = 3
$data[18,11] = 1
$data[18,12] = 10.48767733333333
$data[18,13] = 31.463032
$data[18,14] = 0.1222087640673552
$data[18,15] = 0.1222087640673552
$data[18,16] = 167.3854242795742
$data[18,17] = 1506.468818516168
$data[18,18] = 0.07558416362621881
$data[18,19] = 0.07558416362621881
# row 21: MuSCs -> FAPs
$data[19,0] = "MuSCs"
$data[19,1] = "Ncam1"
$data[19,2] = "Fgfr1"
$data[19,3] = "FAPs"
$data[19,4] = 3
$data[19,5] = 1
$data[19,6] = 15.96019966666667
$data[19,7] = 47.880599
$data[19,8] = 0.6184839868322428
$data[19,9] = 0.6184839868322429
$data[19,10] = 3
$data[19,11] = 1
$data[19,12] = 62.99699166666667
$data[19,13] = 188.990975
$data[19,14] = 0.7340790765058636
$data[19,15] = 0.7340790765058635
$data[19,16] = 1005.444565399336
$data[19,17] = 9049.001088594025
$data[19,18] = 0.4540161538874775
$data[19,19] = 0.4540161538874775
# row 22: MuSCs -> Inflammatory-Mac
$data[20,0] = "MuSCs"
$data[20,1] = "Ncam1"
$data[20,2] = "Fgfr1"
$data[20,3] = "Inflammatory-Mac"
$data[20,4] = 3
$data[20,5] = 1
$data[20,6] = 15.96019966666667
$data[20,7] = 47.880599
$data[20,8] = 0.6184839868322428
$data[20,9] = 0.6184839868322429
$data[20,10] = 3
$data[20,11] = 1
$data[20,12] = 0.3322793333333333
$data[20,13] = 0.996838
$data[20,14] = 0.003871919907635547
$data[20,15] = 0.003871919907635547
$data[20,16] = 5.303244505106889
$data[20,17] = 47.729200545962
$data[20,18] = 0.002394720461169563
$data[20,19] = 0.002394720461169563
# row 23: MuSCs -> MuSCs
$data[21,0] = "MuSCs"
$data[21,1] = "Ncam1"
$data[21,2] = "Fgfr1"
$data[21,3] = "MuSCs"
$data[21,4] = 3
$data[21,5] = 1
$data[21,6] = 15.96019966666667
$data[21,7] = 47.880599
$data[21,8] = 0.6184839868322428
$data[21,9] = 0.6184839868322429
$data[21,10] = 3
$data[21,11] = 1
$data[21,12] = 10.25458433333333
$data[21,13] = 30.763753
$data[21,14] = 0.1194926233493133
$data[21,15] = 0.1194926233493133
$data[21,16] = 163.6652134586719
$data[21,17] = 1472.986921128047
$data[21,18] = 0.07390427408612686
$data[21,19] = 0.07390427408612686
# row 24: MuSCs -> Neutrophils
$data[22,0] = "MuSCs"
$data[22,1] = "Ncam1"
$data[22,2] = "Fgfr1"
$data[22,3] = "Neutrophils"
$data[22,4] = 3
$data[22,5] = 1
$data[22,6] = 15.96019966666667
$data[22,7] = 47.880599
$data[22,8] = 0.6184839868322428
$data[22,9] = 0.6184839868322429
$data[22,10] = 3
$data[22,11] = 1
$data[22,12] = 0.7572163333333334
$data[22,13] = 2.271649
$data[22,14] = 0.008823543029319092
$data[22,15] = 0.00882354302931909
$data[22,16] = 12.08532387086122
$data[22,17] = 108.767914837751
$data[22,18] = 0.005457220070759117
$data[22,19] = 0.005457220070759117
# row 25: MuSCs -> Resolving-Mac
$data[23,0] = "MuSCs"
$data[23,1] = "Ncam1"
$data[23,2] = "Fgfr1"
$data[23,3] = "Resolving-Mac"
$data[23,4] = 3
$data[23,5] = 1
$data[23,6] = 15.96019966666667
$data[23,7] = 47.880599
$data[23,8] = 0.6184839868322428
$data[23,9] = 0.6184839868322429
$data[23,10] = 3
$data[23,11] = 1
$data[23,12] = 0.9889696666666666
$data[23,13] = 2.966909
$data[23,14] = 0.01152407314051338
$data[23,15] = 0.01152407314051338
$data[23,16] = 15.78415334427678
$data[23,17] = 142.057380098491
$data[23,18] = 0.007127454700491079
$data[23,19] = 0.007127454700491079
# row 26: Neutrophils -> ECs
$data[24,0] = "Neutrophils"
$data[24,1] = "Ncam1"
$data[24,2] = "Fgfr1"
$data[24,3] = "ECs"
$data[24,4] = 3
$data[24,5] = 1
$data[24,6] = 0.3268106666666666
$data[24,7] = 0.980432
$data[24,8] = 0.01266445083901121
$data[24,9] = 0.01266445083901121
$data[24,10] = 3
$data[24,11] = 1
$data[24,12] = 10.48767733333333
$data[24,13] = 31.463032
$data[24,14] = 0.1222087640673552
$data[24,15] = 0.1222087640673552
$data[24,16] = 3.427484821091555
$data[24,17] = 30.847363389824
$data[24,18] = 0.00154770688462734
$data[24,19] = 0.00154770688462734
# row 27: Neutrophils -> FAPs
$data[25,0] = "Neutrophils"
$data[25,1] = "Ncam1"
$data[25,2] = "Fgfr1"
$data[25,3] = "FAPs"
$data[25,4] = 3
$data[25,5] = 1
$data[25,6] = 0.3268106666666666
$data[25,7] = 0.980432
$data[25,8] = 0.01266445083901121
$data[25,9] = 0.01266445083901121
$data[25,10] = 3
$data[25,11] = 1
$data[25,12] = 62.99699166666667
$data[25,13] = 188.990975
$data[25,14] = 0.7340790765058636
$data[25,15] = 0.7340790765058635
$data[25,16] = 20.58808884457778
$data[25,17] = 185.2927996012
$data[25,18] = 0.00929670837635526
$data[25,19] = 0.00929670837635526
# row 28: Neutrophils -> Inflammatory-Mac
$data[26,0] = "Neutrophils"
$data[26,1] = "Ncam1"
$data[26,2] = "Fgfr1"
$data[26,3] = "Inflammatory-Mac"
$data[26,4] = 3
$data[26,5] = 1
$data[26,6] = 0.3268106666666666
$data[26,7] = 0.980432
$data[26,8] = 0.01266445083901121
$data[26,9] = 0.01266445083901121
$data[26,10] = 3
$data[26,11] = 1
$data[26,12] = 0.3322793333333333
$data[26,13] = 0.996838
$data[26,14] = 0.003871919907635547
$data[26,15] = 0.003871919907635547
$data[26,16] = 0.1085924304462222
$data[26,17] = 0.977331874016
$data[26,18] = 0.00004903573932283922
$data[26,19] = 0.00004903573932283922
# row 29: Neutrophils -> MuSCs
$data[27,0] = "Neutrophils"
$data[27,1] = "Ncam1"
$data[27,2] = "Fgfr1"
$data[27,3] = "MuSCs"
$data[27,4] = 3
$data[27,5] = 1
$data[27,6] = 0.3268106666666666
$data[27,7] = 0.980432
$data[27,8] = 0.01266445083901121
$data[27,9] = 0.01266445083901121
$data[27,10] = 3
$data[27,11] = 1
$data[27,12] = 10.25458433333333
$data[27,13] = 30.763753
$data[27,14] = 0.1194926233493133
$data[27,15] = 0.1194926233493133
$data[27,16] = 3.351307542366222
$data[27,17] = 30.161767881296
$data[27,18] = 0.001513308454031862
$data[27,19] = 0.001513308454031862
# row 30: Neutrophils -> Neutrophils
$data[28,0] = "Neutrophils"
$data[28,1] = "Ncam1"
$data[28,2] = "Fgfr1"
$data[28,3] = "Neutrophils"
$data[28,4] = 3
$data[28,5] = 1
$data[28,6] = 0.3268106666666666
$data[28,7] = 0.980432
$data[28,8] = 0.01266445083901121
$data[28,9] = 0.01266445083901121
$data[28,10] = 3
$data[28,11] = 1
$data[28,12] = 0.7572163333333334
$data[28,13] = 2.271649
$data[28,14] = 0.008823543029319092
$data[28,15] = 0.00882354302931909
$data[28,16] = 0.2474663747075556
$data[28,17] = 2.227197372368
$data[28,18] = 0.0001117453269207117
$data[28,19] = 0.0001117453269207117
# row 31: Neutrophils -> Resolving-Mac
$data[29,0] = "Neutrophils"
$data[29,1] = "Ncam1"
$data[29,2] = "Fgfr1"
$data[29,3] = "Resolving-Mac"
$data[29,4] = 3
$data[29,5] = 1
$data[29,6] = 0.3268106666666666
$data[29,7] = 0.980432
$data[29,8] = 0.01266445083901121
$data[29,9] = 0.01266445083901121
$data[29,10] = 3
$data[29,11] = 1
$data[29,12] = 0.9889696666666666
$data[29,13] = 2.966909
$data[29,14] = 0.01152407314051338
$data[29,15] = 0.01152407314051338
$data[29,16] = 0.3232058360764444
$data[29,17] = 2.908852524688
$data[29,18] = 0.0001459460577532012
$data[29,19] = 0.0001459460577532012
# row 32: Resolving-Mac -> ECs
$data[30,0] = "Resolving-Mac"
$data[30,1] = "Ncam1"
$data[30,2] = "Fgfr1"
$data[30,3] = "ECs"
$data[30,4] = 1
$data[30,5] = 0.3333333333333333
$data[30,6] = 0.008676333333333333
$data[30,7] = 0.026029
$data[30,8] = 0.0003362221866367304
$data[30,9] = 0.0003362221866367304
$data[30,10] = 3
$data[30,11] = 1
$data[30,12] = 10.48767733333333
$data[30,13] = 31.463032
$data[30,14] = 0.1222087640673552
$data[30,15] = 0.1222087640673552
$data[30,16] = 0.09099458443644443
$data[30,17] = 0.818951259928
$data[30,18] = 0.00004108929788089847
$data[30,19] = 0.00004108929788089847
# row 33: Resolving-Mac -> FAPs
$data[31,0] = "Resolving-Mac"
$data[31,1] = "Ncam1"
$data[31,2] = "Fgfr1"
$data[31,3] = "FAPs"
$data[31,4] = 1
$data[31,5] = 0.3333333333333333
$data[31,6] = 0.008676333333333333
$data[31,7] = 0.026029
$data[31,8] = 0.0003362221866367304
$data[31,9] = 0.0003362221866367304
$data[31,10] = 3
$data[31,11] = 1
$data[31,12] = 62.99699166666667
$data[31,13] = 188.990975
$data[31,14] = 0.7340790765058636
$data[31,15] = 0.7340790765058635
$data[31,16] = 0.5465828986972222
$data[31,17] = 4.919246088275
$data[31,18] = 0.0002468136722670731
$data[31,19] = 0.0002468136722670731
# row 34: Resolving-Mac -> Inflammatory-Mac
$data[32,0] = "Resolving-Mac"
$data[32,1] = "Ncam1"
$data[32,2] = "Fgfr1"
$data[32,3] = "Inflammatory-Mac"
$data[32,4] = 1
$data[32,5] = 0.3333333333333333
$data[32,6] = 0.008676333333333333
$data[32,7] = 0.026029
$data[32,8] = 0.0003362221866367304
$data[32,9] = 0.0003362221866367304
$data[32,10] = 3
$data[32,11] = 1
$data[32,12] = 0.3322793333333333
$data[32,13] = 0.996838
$data[32,14] = 0.003871919907635547
$data[32,15] = 0.003871919907635547
$data[32,16] = 0.002882966255777777
$data[32,17] = 0.025946696302
$data[32,18] = 0.000001301825377827511
$data[32,19] = 0.000001301825377827511
# row 35: Resolving-Mac -> MuSCs
$data[33,0] = "Resolving-Mac"
$data[33,1] = "Ncam1"
$data[33,2] = "Fgfr1"
$data[33,3] = "MuSCs"
$data[33,4] = 1
$data[33,5] = 0.3333333333333333
$data[33,6] = 0.008676333333333333
$data[33,7] = 0.026029
$data[33,8] = 0.0003362221866367304
$data[33,9] = 0.0003362221866367304
$data[33,10] = 3
$data[33,11] = 1
$data[33,12] = 10.25458433333333
$data[33,13] = 30.763753
$data[33,14] = 0.1194926233493133
$data[33,15] = 0.1194926233493133
$data[33,16] = 0.08897219187077778
$data[33,17] = 0.8007497268370001
$data[33,18] = 0.00004017607110946536
$data[33,19] = 0.00004017607110946536
# row 36: Resolving-Mac -> Neutrophils
$data[34,0] = "Resolving-Mac"
$data[34,1] = "Ncam1"
$data[34,2] = "Fgfr1"
$data[34,3] = "Neutrophils"
$data[34,4] = 1
$data[34,5] = 0.3333333333333333
$data[34,6] = 0.008676333333333333
$data[34,7] = 0.026029
$data[34,8] = 0.0003362221866367304
$data[34,9] = 0.0003362221866367304
$data[34,10] = 3
$data[34,11] = 1
$data[34,12] = 0.7572163333333334
$data[34,13] = 2.271649
$data[34,14] = 0.008823543029319092
$data[34,15] = 0.00882354302931909
$data[34,16] = 0.006569861313444445
$data[34,17] = 0.059128751821
$data[34,18] = 0.000002966670931200945
$data[34,19] = 0.000002966670931200945
# row 37: Resolving-Mac -> Resolving-Mac
$data[35,0] = "Resolving-Mac"
$data[35,1] = "Ncam1"
$data[35,2] = "Fgfr1"
$data[35,3] = "Resolving-Mac"
$data[35,4] = 1
$data[35,5] = 0.3333333333333333
$data[35,6] = 0.008676333333333333
$data[35,7] = 0.026029
$data[35,8] = 0.0003362221866367304
$data[35,9] = 0.0003362221866367304
$data[35,10] = 3
$data[35,11] = 1
$data[35,12] = 0.9889696666666666
$data[35,13] = 2.966909
$data[35,14] = 0.01152407314051338
$data[35,15] = 0.01152407314051338
$data[35,16] = 0.008580630484555554
$data[35,17] = 0.07722567436099999
$data[35,18] = 0.000003874649070265021
$data[35,19] = 0.000003874649070265021

$ws.Range("A2:T37").Value = $data
